$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.990.79"
$ws.Range("E2").Value = "  -0.44%  "
$ws.Range("D3").Value = "1.881.11"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'242.81"
$ws.Range("E5").Value = "  -3.84%  "
$ws.Range("D6").Value = "'0.9994"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").Value = "'0.4920"
$ws.Range("E7").Value = "  -3.27%  "
$ws.Range("D8").Value = "'0.2944"
$ws.Range("E8").Value = "  -2.97%  "
$ws.Range("D9").Value = "'0.06625"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").Value = "1.878.34"
$ws.Range("E10").Value = "  -1.54%  "
$ws.Range("D11").Value = "'16.71"
$ws.Range("E11").Value = "  -3.49%  "
$ws.Range("D12").Value = "'0.07174"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "'0.6679"
$ws.Range("E13").Value = "  -3.43%  "
$ws.Range("D14").Value = "'86.91"
$ws.Range("E14").Value = "  +0.00%  "
$ws.Range("D15").Value = "'4.891"
$ws.Range("E15").Value = "  -0.51%  "
$ws.Range("D16").Value = "29.959.67"
$ws.Range("E16").Value = "  -0.52%  "
$ws.Range("D17").Value = "'0.000007825"
$ws.Range("E17").Value = "  -4.48%  "
$ws.Range("D18").Value = "'0.9991"
$ws.Range("E18").Value = "  +0.08%  "
$ws.Range("D19").Value = "'12.79"
$ws.Range("E19").Value = "  -2.25%  "
$ws.Range("D20").Value = "2.120.02"
$ws.Range("E20").Value = "  -1.54%  "
$ws.Range("D21").Value = "'0.9984"
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("D22").Value = "'4.781"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Value = "'5.859"
$ws.Range("E23").Value = "  +1.96%  "
$ws.Range("D24").Value = "'9.110"
$ws.Range("E24").Value = "  -1.92%  "
$ws.Range("D25").Value = "'150.32"
$ws.Range("E25").Value = "  +1.90%  "
$ws.Range("D26").Value = "'140.98"
$ws.Range("E26").Value = "  +4.17%  "
$ws.Range("D27").Value = "'17.01"
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").Value = "'1.908"
$ws.Range("E28").Value = "  -4.68%  "
$ws.Range("D29").Value = "'1.389"
$ws.Range("E29").Value = "  -0.96%  "
$ws.Range("D30").Value = "'4.210"
$ws.Range("E30").Value = "  -1.71%  "
$ws.Range("D31").Value = "'0.08747"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").Value = "'3.985"
$ws.Range("E32").Value = "  -0.63%  "
$ws.Range("D33").Value = "'0.05036"
$ws.Range("E33").Value = "  -0.44%  "
$ws.Range("D34").Value = "'0.7183"
$ws.Range("E34").Value = "  -0.83%  "
$ws.Range("D35").Value = "'1.113"
$ws.Range("E35").Value = "  -2.57%  "
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "'0.01804"
$ws.Range("E37").Value = "  +6.52%  "
$ws.Range("D38").Value = "'2.697"
$ws.Range("E38").Value = "  -4.14%  "
$ws.Range("D39").Value = "'2.159"
$ws.Range("E39").Value = "  -4.98%  "
$ws.Range("D40").Value = "'0.9387"
$ws.Range("E40").Value = "  -2.70%  "
$ws.Range("D41").Value = "'0.4237"
$ws.Range("E41").Value = "  -1.68%  "
$ws.Range("B42").Value = "PaxDollar"
$ws.Range("C42").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D42").Value = "'0.9989"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'103.79"
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'5.745"
$ws.Range("E44").Value = "  -6.64%  "
$ws.Range("D45").Value = "'7.353"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("D46").Value = "'0.1269"
$ws.Range("E46").Value = "  -0.95%  "
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("D48").Value = "'32.64"
$ws.Range("E48").Value = "  -1.62%  "
$ws.Range("D49").Value = "'8.287"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").Value = "'0.3769"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D51").Value = "'56.03"
$ws.Range("E51").Value = "  -1.48%  "
